# Lab Exam 03 grading workbook - fill in grading details for questions
# 18/19 (rows 29-30) and the compilation-error deduction row (37), per the
# "from 33-41 - Driver" grading pass.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 29 - Q18 ("CustomerMappingDriver Class"): partial credit + comment
$ws.Range("E29").Value = 8
$ws.Range("F29").Value = "Partial marks for writing code partially"

# Row 30 - Q19: no points earned, with a comment explaining why
$ws.Range("E30").Value = 0
$ws.Range("F30").Value = "For no output"

# Row 37 - compilation errors deduction: partial deduction + comment
$ws.Range("E37").Value = -2.5
$ws.Range("F37").Value = "For getting exceptions"

# Leave the cursor/selection where grading left off
$ws.Range("F37").Select() | Out-Null
